$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + reporting week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Reference cells used to normalize style indices after type changes ---
# Style 14 = general/text style (used for blank "-" and "***.*" placeholder cells)
# Style 15 = "#,##0" numeric style
# Style 16 = "#,##0.0" numeric style (percent-change columns)
$styleTextRef = $ws.Range("A15")
$styleNumRef15 = $ws.Range("I22")
$styleNumRef16 = $ws.Range("K22")

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = -46.153846153846
$ws.Range("M15").Value = -30
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 72
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -4
$ws.Range("M16").Value = -2.702702702702
$ws.Range("N16").Value = -91.251518833535
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -62.5
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 108
$ws.Range("K17").Value = -2.777777777777
$ws.Range("L17").Value = 12.903225806451
$ws.Range("M17").Value = 15.384615384615
$ws.Range("N17").Value = -66.772151898734
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 80
$ws.Range("J18").Value = 142
$ws.Range("K18").Value = -43.661971830985
$ws.Range("L18").Value = -36
$ws.Range("M18").Value = -33.884297520661
$ws.Range("N18").Value = -93.948562783661
$ws.Range("C19").Value = 37
$ws.Range("D19").Value = 29
$ws.Range("E19").Value = 27.586206896551
$ws.Range("F19").Value = 145
$ws.Range("G19").Value = 164
$ws.Range("H19").Value = -11.585365853658
$ws.Range("I19").Value = 1109
$ws.Range("J19").Value = 1142
$ws.Range("K19").Value = -2.889667250437
$ws.Range("L19").Value = 77.156549520766
$ws.Range("M19").Value = 14.094650205761
$ws.Range("N19").Value = -75.663813912661
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -10
$ws.Range("J20").Value = 74
$ws.Range("K20").Value = -37.837837837837
$ws.Range("M20").Value = 58.620689655172
$ws.Range("N20").Value = -84.768211920529
$ws.Range("C21").Value = 46
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = -2.127659574468
$ws.Range("F21").Value = 180
$ws.Range("G21").Value = 215
$ws.Range("H21").Value = -16.279069767441
$ws.Range("I21").Value = 1421
$ws.Range("J21").Value = 1571
$ws.Range("K21").Value = -9.548058561425
$ws.Range("L21").Value = 46.646026831785
$ws.Range("M21").Value = 9.307692307692
$ws.Range("N21").Value = -80.661404463799
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 6
$ws.Range("I22").Value = 42
$ws.Range("J22").Value = 41
$ws.Range("K22").Value = 2.439024390243
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = 23.529411764705
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = 50
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = -3.508771929824
$ws.Range("F24").Value = 199
$ws.Range("G24").Value = 258
$ws.Range("H24").Value = -22.868217054263
$ws.Range("I24").Value = 1641
$ws.Range("J24").Value = 1545
$ws.Range("K24").Value = 6.213592233009
$ws.Range("L24").Value = 48.238482384823
$ws.Range("M24").Value = 41.465517241379
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 53
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 438
$ws.Range("J25").Value = 341
$ws.Range("K25").Value = 28.445747800586
$ws.Range("L25").Value = 74.501992031872
$ws.Range("M25").Value = 54.770318021201
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = -25
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = -29.166666666666
$ws.Range("L26").Value = -22.727272727272
$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 600
$ws.Range("F27").Value = 18
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 68
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = 19.298245614035
$ws.Range("L27").Value = 30.769230769230

# --- Cells switching from numeric to text placeholder ("0" or "***.*") ---
$c = $ws.Range("F14")
$c.NumberFormat = "@"
$c.Value = "0"
$styleTextRef.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "0"
$styleTextRef.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = "0"
$styleTextRef.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = "0"
$styleTextRef.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = "0"
$styleTextRef.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("H28")
$c.NumberFormat = "@"
$c.Value = "***.*"
$styleTextRef.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = "0"
$styleTextRef.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("H29")
$c.NumberFormat = "@"
$c.Value = "***.*"
$styleTextRef.Copy()
$c.PasteSpecial(-4122)

# --- Cells switching from text placeholder to a real numeric value ---
$c = $ws.Range("C22")
$c.Value = 1
$styleNumRef15.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("D22")
$c.Value = 2
$styleNumRef15.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("E22")
$c.Value = -50
$styleNumRef16.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("D23")
$c.Value = 1
$styleNumRef15.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("E23")
$c.Value = -100
$styleNumRef16.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("G23")
$c.Value = 1
$styleNumRef15.Copy()
$c.PasteSpecial(-4122)
$c = $ws.Range("H23")
$c.Value = -100
$styleNumRef16.Copy()
$c.PasteSpecial(-4122)

$excel.CutCopyMode = 0
